# Shaconemo ocean/ocnBgchem/seaIce ping-file revision bump: r270 -> r274.
# These two Shaconemo revisions are equivalent for ece2cmor3 purposes, so
# the only functional change is the revision number quoted inside the
# "comment" column (column H) for every data row (rows 3 through 282).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$commentRange = $ws.Range("H3:H282")
$commentRange.Replace("r270", "r274")

# Leave the sheet scrolled back to the top with the edited comment
# column selected, matching the view state saved with the workbook.
$commentRange.Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
